$wb = $excel.ActiveWorkbook

# --- Append new Week 13 play-log entries to the running shared-string logs on the YDS sheet ---
$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value2 = $wsYDS.Range("B2").Value2 + ' 0 2 2 4 11 2 0 6 1 6 2 3 0 2 3 4 1 3 4 4 2 3 2 2'
$wsYDS.Range("B3").Value2 = $wsYDS.Range("B3").Value2 + ' 9 5 11 1 4 12 8 6 -2 6 8 4 8 5 9 2 9 17 4 6 5 3 17 25 11 13 3 2 16 17'
$wsYDS.Range("C2").Value2 = $wsYDS.Range("C2").Value2 + ' 3 3 10 -2 0 1 0 7 0 23 17 1 0 9 13 3'
$wsYDS.Range("C3").Value2 = $wsYDS.Range("C3").Value2 + ' 3 15 -3 7 -5 4 18 20 -1 5 12 16 6 11 17 8 3 9 8 2 10 11 11'

# --- Append new Week 13 play-log entries to the running shared-string logs on the ST sheet ---
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B6").Value2 = $wsST.Range("B6").Value2 + ' 15 17'
$wsST.Range("D3").Value2 = $wsST.Range("D3").Value2 + ' 53 50 65 39 37 46'
$wsST.Range("D4").Value2 = $wsST.Range("D4").Value2 + ' 0 15 11 0 -4 0'
$wsST.Range("D5").Value2 = $wsST.Range("D5").Value2 + ' 0 0 7 0 0 0'

# --- Update running season totals for each stat sheet ---

$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value2 = 321
$ws.Range("E2").Value2 = 18
$ws.Range("F2").Value2 = 110
$ws.Range("G2").Value2 = 100
$ws.Range("J2").Value2 = 57
$ws.Range("L2").Value2 = 548
$ws.Range("M2").Value2 = 362
$ws.Range("Q2").Value2 = 971
$ws.Range("C3").Value2 = 366
$ws.Range("E3").Value2 = 45
$ws.Range("F3").Value2 = 194
$ws.Range("G3").Value2 = 87
$ws.Range("H3").Value2 = 47
$ws.Range("I3").Value2 = 93
$ws.Range("J3").Value2 = 121
$ws.Range("N3").Value2 = 35

$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value2 = 345
$ws.Range("F2").Value2 = 99
$ws.Range("G2").Value2 = 85
$ws.Range("H2").Value2 = 6
$ws.Range("I2").Value2 = 12
$ws.Range("L2").Value2 = 553
$ws.Range("M2").Value2 = 330
$ws.Range("O2").Value2 = 41
$ws.Range("P2").Value2 = 21
$ws.Range("Q2").Value2 = 936
$ws.Range("C3").Value2 = 356
$ws.Range("D3").Value2 = 6
$ws.Range("E3").Value2 = 63
$ws.Range("F3").Value2 = 212
$ws.Range("H3").Value2 = 59
$ws.Range("I3").Value2 = 109
$ws.Range("J3").Value2 = 94
$ws.Range("N3").Value2 = 44

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value2 = 144
$ws.Range("D2").Value2 = 127
$ws.Range("F2").Value2 = 128
$ws.Range("G2").Value2 = 125
$ws.Range("L2").Value2 = 40
$ws.Range("M2").Value2 = 30
$ws.Range("N2").Value2 = 19
$ws.Range("B3").Value2 = 104

$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("C2").Value2 = 20
$ws.Range("D3").Value2 = 14

$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value2 = 31
$ws.Range("D2").Value2 = 11
